# Alteração nos rótulos da tabela para já transformar a primeira linha em
# cabeçalho automaticamente no Power BI.
#
# Sheets 1,2,3,5 (and sheet 6, only column B) use year labels (2015, 2030,
# 2040, 2050) in row 1 and get prefixed with "Ano ".
# Sheet 4 uses interval labels (2015, 2015-2030, 2031-2040, 2041-2050) in
# row 1 and gets prefixed with "Intervalo ".

$wb = $excel.ActiveWorkbook

$anoSheets = @(
    "Potencia Acumulada - SIN (MW)",
    "Geracao Periodo Medio (MWMed)",
    "Atendimento a Ponta(MW)",
    "Emissoes Totais (MtCO2eq)"
)

$intervaloSheets = @(
    "Potencia Incremental - SIN(MW)"
)

foreach ($name in $anoSheets) {
    $ws = $wb.Worksheets.Item($name)
    foreach ($col in @("B", "C", "D", "E")) {
        $cell = $ws.Range($col + "1")
        $current = $cell.Text
        $cell.Value = "Ano " + $current
    }
}

foreach ($name in $intervaloSheets) {
    $ws = $wb.Worksheets.Item($name)
    foreach ($col in @("B", "C", "D", "E")) {
        $cell = $ws.Range($col + "1")
        $current = $cell.Text
        $cell.Value = "Intervalo " + $current
    }
}

$wsCusto = $wb.Worksheets.Item("Custo Total (bilhões de R$)")
$cell = $wsCusto.Range("B1")
$current = $cell.Text
$cell.Value = "Ano " + $current
